$p = $ppt.ActivePresentation

# Delete slides that are removed in the target deck (delete from the
# highest index down so earlier indices stay valid):
#   slide 7 -> old sldId 407 "Individual Presentations of Research Projects"
#   slide 6 -> old sldId 409 "Upcoming Classes (until the end of semester)"
#   slide 3 -> old sldId 412 "Proposal Abstract Review Criteria"
#   slide 2 -> old sldId 411 "Proposal Abstract Review Panel (In-Class Activity)"
$p.Slides.Item(7).Delete()
$p.Slides.Item(6).Delete()
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()

# Remaining slide 1 (title slide): bump the class number and drop the
# "Proposal Abstract Review Panel" bullet line (that topic's slides were
# just removed above), keeping the "Physics Outreach" bullet intact.
$s1 = $p.Slides.Item(1)
$shp = $s1.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Remove the whole "Proposal Abstract Review Panel" paragraph (2nd
# paragraph), including its trailing paragraph mark, so the following
# "Physics Outreach" paragraph keeps its own bullet formatting.
$tr.Paragraphs(2, 1).Delete()

# Update "Class #24" -> "Class #23" in place (select just the text run,
# not the trailing paragraph mark, so it stays a single run).
$tr.Characters(1, 9).Text = "Class #23"
